$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gm13306"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1181316666666667
$ws.Range("H2").Value = 0.354395
$ws.Range("I2").Value = 0.1921951611040097
$ws.Range("J2").Value = 0.1921951611040097
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.5706193333333334
$ws.Range("N2").Value = 1.711858
$ws.Range("O2").Value = 0.4188640502130462
$ws.Range("P2").Value = 0.4188640502130463
$ws.Range("Q2").Value = 0.0674082128788889
$ws.Range("R2").Value = 0.6066739159100001
$ws.Range("S2").Value = 0.08050364361137445
$ws.Range("T2").Value = 0.08050364361137446

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gm13306"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1181316666666667
$ws.Range("H3").Value = 0.354395
$ws.Range("I3").Value = 0.1921951611040097
$ws.Range("J3").Value = 0.1921951611040097
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4846943333333333
$ws.Range("N3").Value = 1.454083
$ws.Range("O3").Value = 0.3557906641356566
$ws.Range("P3").Value = 0.3557906641356566
$ws.Range("Q3").Value = 0.05725774942055555
$ws.Range("R3").Value = 0.515319744785
$ws.Range("S3").Value = 0.06838124401285514
$ws.Range("T3").Value = 0.06838124401285516

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gm13306"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1181316666666667
$ws.Range("H4").Value = 0.354395
$ws.Range("I4").Value = 0.1921951611040097
$ws.Range("J4").Value = 0.1921951611040097
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3069883333333334
$ws.Range("N4").Value = 0.920965
$ws.Range("O4").Value = 0.2253452856512971
$ws.Range("P4").Value = 0.2253452856512971
$ws.Range("Q4").Value = 0.0362650434638889
$ws.Range("R4").Value = 0.326385391175
$ws.Range("S4").Value = 0.04331027347978014
$ws.Range("T4").Value = 0.04331027347978014

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gm13306"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4965126666666667
$ws.Range("H5").Value = 1.489538
$ws.Range("I5").Value = 0.8078048388959902
$ws.Range("J5").Value = 0.8078048388959903
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.5706193333333334
$ws.Range("N5").Value = 1.711858
$ws.Range("O5").Value = 0.4188640502130462
$ws.Range("P5").Value = 0.4188640502130463
$ws.Range("Q5").Value = 0.2833197268448889
$ws.Range("R5").Value = 2.549877541604
$ws.Range("S5").Value = 0.3383604066016718
$ws.Range("T5").Value = 0.3383604066016718

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gm13306"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4965126666666667
$ws.Range("H6").Value = 1.489538
$ws.Range("I6").Value = 0.8078048388959902
$ws.Range("J6").Value = 0.8078048388959903
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4846943333333333
$ws.Range("N6").Value = 1.454083
$ws.Range("O6").Value = 0.3557906641356566
$ws.Range("P6").Value = 0.3557906641356566
$ws.Range("Q6").Value = 0.2406568759615555
$ws.Range("R6").Value = 2.165911883654
$ws.Range("S6").Value = 0.2874094201228014
$ws.Range("T6").Value = 0.2874094201228015

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gm13306"
$ws.Range("C7").Value = "Ccr10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4965126666666667
$ws.Range("H7").Value = 1.489538
$ws.Range("I7").Value = 0.8078048388959902
$ws.Range("J7").Value = 0.8078048388959903
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3069883333333334
$ws.Range("N7").Value = 0.920965
$ws.Range("O7").Value = 0.2253452856512971
$ws.Range("P7").Value = 0.2253452856512971
$ws.Range("Q7").Value = 0.1524235960188889
$ws.Range("R7").Value = 1.37181236417
$ws.Range("S7").Value = 0.182035012171517
$ws.Range("T7").Value = 0.182035012171517

